# Weekly data refresh: insert this week's new price record for
# "Femacal de La Calera - Ají" ahead of the historical rows, shifting the
# existing rows (old 256..301) down by one (new 257..302).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 256; Excel shifts every
# row at/after 256 down by one (old 256 -> 257, ..., old 301 -> 302) and
# the sheet's used range / dimension grows from R301 to R302 automatically.
$ws.Rows.Item(256).EntireRow.Insert()

# Populate the newly inserted row 256 with this week's record.
$ws.Range("A256").Value = 3
$ws.Range("B256").Value = "Femacal de La Calera"
$ws.Range("C256").Value = "Coquimbo"
$ws.Range("D256").Value = 44505
$ws.Range("E256").Value = 5
$ws.Range("F256").Value = 100112021
$ws.Range("G256").Value = "Ají"
$ws.Range("H256").Value = "Americana (o)"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 67
$ws.Range("K256").Value = 41000
$ws.Range("L256").Value = 42000
$ws.Range("M256").Value = 41478
$ws.Range("N256").Value = "$/caja 15 kilos"
$ws.Range("O256").Value = "Región de Arica y Parinacota"
$ws.Range("P256").Value = 2765
$ws.Range("Q256").Value = 15
$ws.Range("R256").Value = "Hortaliza"
